$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 15 (G=44146)
$ws.Range("H15").Value = 1014.5455
$ws.Range("I15").Value = 1014.5455
$ws.Range("K15").Value = 3043.6365
$ws.Range("M15").Value = -2874.6365
# Row 28 (G=27772)
$ws.Range("H28").Value = 750.75
$ws.Range("I28").Value = 658.4211
$ws.Range("J28").Value = 2505
$ws.Range("K28").Value = 658.4211
$ws.Range("L28").Value = 2505
$ws.Range("M28").Value = -173.4211
$ws.Range("N28").Value = -3475
# Row 112 (G=27960)
$ws.Range("H112").Value = 912627.6
$ws.Range("J112").Value = 1254669.4
$ws.Range("L112").Value = 3764008.2
$ws.Range("N112").Value = -3766224.2
# Row 116 (G=27778)
$ws.Range("H116").Value = 3534.3845
$ws.Range("I116").Value = 3578.9167
$ws.Range("J116").Value = 3000
$ws.Range("K116").Value = 3578.9167
$ws.Range("L116").Value = 3000
$ws.Range("M116").Value = -136.9167000000002
$ws.Range("N116").Value = -9884
# Row 129 (G=36115)
$ws.Range("H129").Value = 111112550
$ws.Range("J129").Value = 2832.3333
$ws.Range("L129").Value = 8496.999899999999
$ws.Range("N129").Value = -18496.9999
# Row 132 (G=44049)
$ws.Range("H132").Value = 8391.548000000001
$ws.Range("I132").Value = 3071.7942
$ws.Range("K132").Value = 9215.382599999999
$ws.Range("M132").Value = -6685.382599999999
# Row 135 (G=44047)
$ws.Range("H135").Value = 13539.667
$ws.Range("I135").Value = 10000
$ws.Range("K135").Value = 90000
$ws.Range("M135").Value = -87465
# Row 138 (G=44169)
$ws.Range("H138").Value = 367881.34
$ws.Range("J138").Value = 1251859.9
$ws.Range("L138").Value = 3755579.7
$ws.Range("N138").Value = -3765859.7

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2 (G=27713)
$ws.Range("H2").Value = 1780.5333
$ws.Range("I2").Value = 1906.3636
$ws.Range("J2").Value = 1434.5
$ws.Range("K2").Value = 1906.3636
$ws.Range("L2").Value = 1434.5
$ws.Range("M2").Value = -1793.3636
$ws.Range("N2").Value = -1660.5
# Row 32 (G=44147)
$ws.Range("H32").Value = 6774.3896
$ws.Range("I32").Value = 6869.982
$ws.Range("K32").Value = 6869.982
$ws.Range("M32").Value = -6582.982
# Row 45 (G=27714)
$ws.Range("H45").Value = 3205.2144
$ws.Range("I45").Value = 2757.5
$ws.Range("K45").Value = 2757.5
$ws.Range("M45").Value = -2380.5
# Row 116 (G=27713)
$ws.Range("H116").Value = 1780.5333
$ws.Range("I116").Value = 1906.3636
$ws.Range("J116").Value = 1434.5
$ws.Range("K116").Value = 1906.3636
$ws.Range("L116").Value = 1434.5
$ws.Range("M116").Value = 387.6364000000001
$ws.Range("N116").Value = -6022.5
# Row 123 (G=34107)
$ws.Range("H123").Value = 59166.168
$ws.Range("J123").Value = 59166.168
$ws.Range("L123").Value = 59166.168
$ws.Range("N123").Value = -68966.16800000001

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3 (G=27713)
$ws.Range("H3").Value = 1780.5333
$ws.Range("I3").Value = 1906.3636
$ws.Range("J3").Value = 1434.5
$ws.Range("K3").Value = 1906.3636
$ws.Range("L3").Value = 1434.5
$ws.Range("M3").Value = -1792.3636
$ws.Range("N3").Value = -1662.5
# Row 98 (G=19545)
$ws.Range("H98").Value = 89999
$ws.Range("J98").Value = 89999
$ws.Range("L98").Value = 89999
$ws.Range("N98").Value = -95989
# Row 134 (G=43998)
$ws.Range("H134").Value = 5253.933
$ws.Range("I134").Value = 1188.5834
$ws.Range("K134").Value = 3565.7502
$ws.Range("M134").Value = -1030.7502

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 22 (G=5367)
$ws.Range("H22").Value = 498.05
$ws.Range("I22").Value = 466.23077
$ws.Range("J22").Value = 557.1429000000001
$ws.Range("K22").Value = 466.23077
$ws.Range("L22").Value = 557.1429000000001
$ws.Range("M22").Value = -116.23077
$ws.Range("N22").Value = -1257.1429
# Row 31 (G=44023)
$ws.Range("H31").Value = 2600.4814
$ws.Range("I31").Value = 1441.591
$ws.Range("K31").Value = 1441.591
$ws.Range("M31").Value = -1146.591
# Row 34 (G=44023)
$ws.Range("H34").Value = 2600.4814
$ws.Range("I34").Value = 1441.591
$ws.Range("K34").Value = 1441.591
$ws.Range("M34").Value = -1239.591
# Row 58 (G=44021)
$ws.Range("H58").Value = 3990.6667
$ws.Range("I58").Value = 5972.5
$ws.Range("K58").Value = 5972.5
$ws.Range("M58").Value = -5769.5
# Row 95 (G=18192)
$ws.Range("H95").Value = 42000
$ws.Range("J95").Value = 42000
$ws.Range("L95").Value = 42000
$ws.Range("N95").Value = -47492
# Row 105 (G=19928)
$ws.Range("H105").Value = 1933.5714
$ws.Range("I105").Value = 1845
$ws.Range("K105").Value = 1845
$ws.Range("M105").Value = -98
# Row 132 (G=44019)
$ws.Range("H132").Value = 1483892.9
$ws.Range("I132").Value = 1820559.6
$ws.Range("K132").Value = 5461678.800000001
$ws.Range("M132").Value = -5459148.800000001
# Row 134 (G=44020)
$ws.Range("H134").Value = 2271.4644
$ws.Range("I134").Value = 1209.159
$ws.Range("J134").Value = 6166.5835
$ws.Range("K134").Value = 3627.477
$ws.Range("L134").Value = 18499.7505
$ws.Range("M134").Value = -1092.477
$ws.Range("N134").Value = -23569.7505
# Row 135 (G=42008)
$ws.Range("H135").Value = 91333.336
$ws.Range("J135").Value = 85000
$ws.Range("L135").Value = 85000
$ws.Range("N135").Value = -95140
# Row 136 (G=44021)
$ws.Range("H136").Value = 3990.6667
$ws.Range("I136").Value = 5972.5
$ws.Range("K136").Value = 17917.5
$ws.Range("M136").Value = -15367.5

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5 (G=43974)
$ws.Range("H5").Value = 1407.8636
$ws.Range("I5").Value = 849
$ws.Range("K5").Value = 2547
$ws.Range("M5").Value = -2435
# Row 37 (G=9516)
$ws.Range("H37").Value = 199995
$ws.Range("J37").Value = 199995
$ws.Range("L37").Value = 599985
$ws.Range("N37").Value = -600209
# Row 122 (G=36078)
$ws.Range("H122").Value = 754.35297
$ws.Range("J122").Value = 804.93335
$ws.Range("L122").Value = 7244.40015
$ws.Range("N122").Value = -12144.40015
# Row 132 (G=43972)
$ws.Range("H132").Value = 1135.7142
$ws.Range("I132").Value = 1123
$ws.Range("J132").Value = 1145.25
$ws.Range("K132").Value = 10107
$ws.Range("L132").Value = 10307.25
$ws.Range("M132").Value = -7577
$ws.Range("N132").Value = -15367.25
# Row 135 (G=43974)
$ws.Range("H135").Value = 1407.8636
$ws.Range("I135").Value = 849
$ws.Range("K135").Value = 7641
$ws.Range("M135").Value = -5106

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 52 (G=4147)
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").Value = ""
# Row 55 (G=4237)
$ws.Range("H55").Value = 10978.625
$ws.Range("I55").Value = 15082.5
$ws.Range("J55").Value = 6874.75
$ws.Range("K55").Value = 15082.5
$ws.Range("L55").Value = 6874.75
$ws.Range("M55").Value = -14755.5
$ws.Range("N55").Value = -7528.75
# Row 102 (G=36169)
$ws.Range("H102").Value = 34562.227
$ws.Range("I102").Value = 2146.0454
$ws.Range("J102").Value = 113801.78
$ws.Range("K102").Value = 2146.0454
$ws.Range("L102").Value = 113801.78
$ws.Range("M102").Value = -524.0454
$ws.Range("N102").Value = -117045.78
# Row 107 (G=27802)
$ws.Range("H107").Value = 731.4400000000001
$ws.Range("I107").Value = 424.2857
$ws.Range("J107").Value = 1122.3636
$ws.Range("K107").Value = 424.2857
$ws.Range("L107").Value = 1122.3636
$ws.Range("M107").Value = 1495.7143
$ws.Range("N107").Value = -4962.3636
# Row 113 (G=27710)
$ws.Range("H113").Value = 1623.9333
$ws.Range("I113").Value = 1550.909
$ws.Range("K113").Value = 1550.909
$ws.Range("M113").Value = 619.0909999999999
# Row 132 (G=44008)
$ws.Range("H132").Value = 13336401
$ws.Range("I132").Value = 13892001
$ws.Range("K132").Value = 41676003
$ws.Range("M132").Value = -41673473
# Row 133 (G=41854)
$ws.Range("H133").Value = 70776.664
$ws.Range("J133").Value = 70776.664
$ws.Range("L133").Value = 70776.664
$ws.Range("N133").Value = -80896.664
# Row 140 (G=42458)
$ws.Range("H140").Value = 80999.664
$ws.Range("J140").Value = 80999.664
$ws.Range("L140").Value = 80999.664
$ws.Range("N140").Value = -91359.664

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 46 (G=5282)
$ws.Range("H46").Value = 8094.1577
$ws.Range("I46").Value = 3749.8333
$ws.Range("K46").Value = 3749.8333
$ws.Range("M46").Value = -3561.8333
# Row 61 (G=27740)
$ws.Range("H61").Value = 6138.615
$ws.Range("I61").Value = 7259.2
$ws.Range("K61").Value = 7259.2
$ws.Range("M61").Value = -7057.2
# Row 100 (G=19995)
$ws.Range("H100").Value = 4107
$ws.Range("I100").Value = 3357.1428
$ws.Range("J100").Value = 4856.857
$ws.Range("K100").Value = 3357.1428
$ws.Range("L100").Value = 4856.857
$ws.Range("M100").Value = -2816.1428
$ws.Range("N100").Value = -5938.857
# Row 108 (G=25655)
$ws.Range("H108").Value = 71198
$ws.Range("J108").Value = 76497.5
$ws.Range("L108").Value = 76497.5
$ws.Range("N108").Value = -84177.5
# Row 113 (G=27740)
$ws.Range("H113").Value = 6138.615
$ws.Range("I113").Value = 7259.2
$ws.Range("K113").Value = 7259.2
$ws.Range("M113").Value = -5089.2
# Row 132 (G=44058)
$ws.Range("H132").Value = 3163.7856
$ws.Range("I132").Value = 3171.6128
$ws.Range("J132").Value = 3141.7273
$ws.Range("K132").Value = 9514.838400000001
$ws.Range("L132").Value = 9425.1819
$ws.Range("M132").Value = -6984.838400000001
$ws.Range("N132").Value = -14485.1819

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 132 (G=44029)
$ws.Range("H132").Value = 2078.149
$ws.Range("I132").Value = 2053.2432
$ws.Range("K132").Value = 6159.7296
$ws.Range("M132").Value = -3629.7296

